$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "34.189.66"
$ws.Range("E2").Value = "  +0.62%  "

$ws.Range("D3").Value = "1.788.48"
$ws.Range("E3").Value = "  +0.53%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "226.69"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.40%  "

$ws.Range("E6").Value = "  -0.98%  "

$ws.Range("E7").Value = "  +0.06%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "31.94"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.90%  "

$ws.Range("E9").Value = "  +0.89%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0689"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.57%  "

$ws.Range("E11").Value = "  +0.70%  "

$ws.Range("D12").Value = "2.045.81"
$ws.Range("E12").Value = "  +0.46%  "

$ws.Range("B13").Value = "Chainlink"
$ws.Range("C13").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "11.03"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.11%  "

$ws.Range("B14").Value = "WrappedEther"
$ws.Range("C14").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D14").Value = "1.779.19"
$ws.Range("E14").Value = "  +0.16%  "

$ws.Range("D15").Value = "34.144.95"
$ws.Range("E15").Value = "  +0.52%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.623"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.46%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "4.18"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.40%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "68.30"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.04%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "246.42"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.66%  "

$ws.Range("E20").Value = "  -0.62%  "

$ws.Range("E21").Value = "  +0.12%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "10.82"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.27%  "

$ws.Range("E23").Value = "  -0.08%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.05"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.11%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "161.16"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.98%  "

$ws.Range("E26").Value = "  +1.16%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "16.35"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.62%  "

$ws.Range("E28").Value = "  +0.69%  "

$ws.Range("E29").Value = "  +0.14%  "

$ws.Range("E30").Value = "  -0.19%  "

$ws.Range("E31").Value = "  +1.16%  "

$ws.Range("E32").Value = "  +0.53%  "

$ws.Range("E33").Value = "  +1.91%  "

$ws.Range("E34").Value = "  -0.34%  "

$ws.Range("D35").Value = "1.442.88"
$ws.Range("E35").Value = "  +3.67%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.649"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.75%  "

$ws.Range("E37").Value = "  +7.62%  "

$ws.Range("E38").Value = "  +2.88%  "

$ws.Range("E39").Value = "  +0.10%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "80.55"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +3.37%  "

$ws.Range("E41").Value = "  +0.54%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.920"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.09%  "

$ws.Range("E43").Value = "  +1.33%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "13.49"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.50%  "

$ws.Range("E45").Value = "  +2.15%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "6.05"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +3.36%  "

$ws.Range("E47").Value = "  -0.48%  "

$ws.Range("D48").Value = "0.0₆0135"
$ws.Range("E48").Value = "  -6.53%  "

$ws.Range("D49").Value = "1.947.72"
$ws.Range("E49").Value = "  +0.69%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "105.83"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.46%  "

$ws.Range("E51").Value = "  +0.11%  "
